$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Crew Stock Report")

$ws.Range("C4").Value = "From : 01-10-2018"
$ws.Range("E4").Value = "To : 30-04-2019"
$ws.Range("F11").Value = "Print taken at : 30-04-2019 14:31:48"
